# chore: update Sheets via scheduled runner
# Applies refreshed profit-model figures (columns H-N) across ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2468243.5
$ws.Range("J17").Value = 2506807.2
$ws.Range("L17").Value = 7520421.600000001
$ws.Range("N17").Value = -7520757.600000001
$ws.Range("H88").Value = 2077.182
$ws.Range("I88").Value = 1650
$ws.Range("J88").Value = 2321.2856
$ws.Range("K88").Value = 1650
$ws.Range("L88").Value = 2321.2856
$ws.Range("M88").Value = -1244
$ws.Range("N88").Value = -3133.2856
$ws.Range("H91").Value = 2077.182
$ws.Range("I91").Value = 1650
$ws.Range("J91").Value = 2321.2856
$ws.Range("K91").Value = 1650
$ws.Range("L91").Value = 2321.2856
$ws.Range("M91").Value = -246
$ws.Range("N91").Value = -5129.2856
$ws.Range("H112").Value = 1586.7755
$ws.Range("I112").Value = 1133.3334
$ws.Range("J112").Value = 1616.3478
$ws.Range("K112").Value = 3400.0002
$ws.Range("L112").Value = 4849.0434
$ws.Range("M112").Value = -2292.0002
$ws.Range("N112").Value = -7065.0434
$ws.Range("H129").Value = 887.62
$ws.Range("I129").Value = 491.44446
$ws.Range("J129").Value = 926.8022
$ws.Range("K129").Value = 1474.33338
$ws.Range("L129").Value = 2780.4066
$ws.Range("M129").Value = 3525.66662
$ws.Range("N129").Value = -12780.4066
$ws.Range("H137").Value = 1429.96
$ws.Range("I137").Value = 1195.7142
$ws.Range("J137").Value = 2659.75
$ws.Range("K137").Value = 3587.1426
$ws.Range("L137").Value = 7979.25
$ws.Range("M137").Value = -1037.1426
$ws.Range("N137").Value = -13079.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4561.6787
$ws.Range("I32").Value = 3666.8
$ws.Range("J32").Value = 8222.546
$ws.Range("K32").Value = 3666.8
$ws.Range("L32").Value = 8222.546
$ws.Range("M32").Value = -3379.8
$ws.Range("N32").Value = -8796.546
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7598.3706
$ws.Range("I31").Value = 2111.5454
$ws.Range("J31").Value = 11370.5625
$ws.Range("K31").Value = 2111.5454
$ws.Range("L31").Value = 11370.5625
$ws.Range("M31").Value = -1816.5454
$ws.Range("N31").Value = -11960.5625
$ws.Range("H34").Value = 7598.3706
$ws.Range("I34").Value = 2111.5454
$ws.Range("J34").Value = 11370.5625
$ws.Range("K34").Value = 2111.5454
$ws.Range("L34").Value = 11370.5625
$ws.Range("M34").Value = -1909.5454
$ws.Range("N34").Value = -11774.5625
$ws.Range("H58").Value = 1490.5186
$ws.Range("I58").Value = 1164.375
$ws.Range("J58").Value = 1964.909
$ws.Range("K58").Value = 1164.375
$ws.Range("L58").Value = 1964.909
$ws.Range("M58").Value = -961.375
$ws.Range("N58").Value = -2370.909
$ws.Range("H134").Value = 3366.84
$ws.Range("I134").Value = 4206.4443
$ws.Range("J134").Value = 1207.8572
$ws.Range("K134").Value = 12619.3329
$ws.Range("L134").Value = 3623.5716
$ws.Range("M134").Value = -10084.3329
$ws.Range("N134").Value = -8693.5716
$ws.Range("H136").Value = 1490.5186
$ws.Range("I136").Value = 1164.375
$ws.Range("J136").Value = 1964.909
$ws.Range("K136").Value = 3493.125
$ws.Range("L136").Value = 5894.727000000001
$ws.Range("M136").Value = -943.125
$ws.Range("N136").Value = -10994.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 8333710.5
$ws.Range("J98").Value = 8333710.5
$ws.Range("L98").Value = 25001131.5
$ws.Range("N98").Value = -25004127.5
$ws.Range("H113").Value = 294615.25
$ws.Range("I113").Value = 495.96155
$ws.Range("J113").Value = 1250503
$ws.Range("K113").Value = 1487.88465
$ws.Range("L113").Value = 3751509
$ws.Range("M113").Value = 682.11535
$ws.Range("N113").Value = -3755849
$ws.Range("H136").Value = 4608.548
$ws.Range("I136").Value = 20591.8
$ws.Range("J136").Value = 2448.6487
$ws.Range("K136").Value = 61775.39999999999
$ws.Range("L136").Value = 7345.946100000001
$ws.Range("M136").Value = -56675.39999999999
$ws.Range("N136").Value = -17545.9461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24800
$ws.Range("J52").Value = 24800
$ws.Range("L52").Value = 24800
$ws.Range("N52").Value = -25318
$ws.Range("H57").Value = 8924.786
$ws.Range("J57").Value = 8924.786
$ws.Range("L57").Value = 8924.786
$ws.Range("N57").Value = -10564.786
$ws.Range("H113").Value = 90910400
$ws.Range("I113").Value = 166667380
$ws.Range("J113").Value = 2020
$ws.Range("K113").Value = 166667380
$ws.Range("L113").Value = 2020
$ws.Range("M113").Value = -166665210
$ws.Range("N113").Value = -6360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2779693.2
$ws.Range("I22").Value = 22222672
$ws.Range("J22").Value = 2124.7715
$ws.Range("K22").Value = 22222672
$ws.Range("L22").Value = 2124.7715
$ws.Range("M22").Value = -22222377
$ws.Range("N22").Value = -2714.7715
$ws.Range("H27").Value = 2779693.2
$ws.Range("I27").Value = 22222672
$ws.Range("J27").Value = 2124.7715
$ws.Range("K27").Value = 22222672
$ws.Range("L27").Value = 2124.7715
$ws.Range("M27").Value = -22222565
$ws.Range("N27").Value = -2338.7715
$ws.Range("H132").Value = 24130912
$ws.Range("I132").Value = 28654472
$ws.Range("J132").Value = 5266
$ws.Range("K132").Value = 85963416
$ws.Range("L132").Value = 15798
$ws.Range("M132").Value = -85960886
$ws.Range("N132").Value = -20858
$ws.Range("H136").Value = 6104.9375
$ws.Range("I136").Value = 7161
$ws.Range("J136").Value = 3406.111
$ws.Range("K136").Value = 21483
$ws.Range("L136").Value = 10218.333
$ws.Range("M136").Value = -18933
$ws.Range("N136").Value = -15318.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 41600
$ws.Range("J68").Value = 41600
$ws.Range("L68").Value = 41600
$ws.Range("N68").Value = -43222
$ws.Range("H69").Value = 30135.5
$ws.Range("J69").Value = 30135.5
$ws.Range("L69").Value = 30135.5
$ws.Range("N69").Value = -31633.5
$ws.Range("H71").Value = 41600
$ws.Range("J71").Value = 41600
$ws.Range("L71").Value = 124800
$ws.Range("N71").Value = -132912
$ws.Range("H72").Value = 30135.5
$ws.Range("J72").Value = 30135.5
$ws.Range("L72").Value = 90406.5
$ws.Range("N72").Value = -97894.5
$ws.Range("H76").Value = 42200
$ws.Range("J76").Value = 42200
$ws.Range("L76").Value = 42200
$ws.Range("N76").Value = -42830
$ws.Range("H79").Value = 42200
$ws.Range("J79").Value = 42200
$ws.Range("L79").Value = 42200
$ws.Range("N79").Value = -44384
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H122").Value = 3007.158
$ws.Range("I122").Value = 2503
$ws.Range("J122").Value = 4099.5
$ws.Range("K122").Value = 7509
$ws.Range("L122").Value = 12298.5
$ws.Range("M122").Value = -5059
$ws.Range("N122").Value = -17198.5
